# Update the "想去人数" (F column) figures that changed between the two
# data snapshots. The same values need to be updated on both the "展览"
# sheet and the "全部类型" sheet, since they mirror the same rows.

$wb = $excel.ActiveWorkbook

# Map of row number -> new F-column value (only rows whose value changed)
$updates = @{
    2  = 265
    3  = 1366
    10 = 134
    11 = 4591
    12 = 6859
    16 = 572
    17 = 55
    18 = 4137
    19 = 559
    20 = 76
    21 = 61
    22 = 2717
    26 = 360
    27 = 369
    31 = 1630
    32 = 1023
    33 = 66
    34 = 197
    35 = 84
    40 = 134
    41 = 645
    42 = 15
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
